$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("E2").Value = 23.20000000000019
$ws.Range("G2").Value = 0.0002422121778209441
$ws.Range("H2").Value = 0.003059608378842459
$ws.Range("I2").Value = 0.7142584581539119
$ws.Range("K2").Value = 4.791519168471507
$ws.Range("L2").Value = "[1.6599565034869075, 7.923081833456107]"
$ws.Range("M2").Value = 0.002790137072895771
$ws.Range("N2").Value = 0.002790137072895771
$ws.Range("O2").Value = -0.729579074707539
$ws.Range("P2").Value = "[-1.3333686537758478, -0.12578949563923025]"
$ws.Range("Q2").Value = 0.01798819731824608
$ws.Range("R2").Value = 0.03597639463649216
$ws.Range("S2").Value = 14.31080508464335
$ws.Range("T2").Value = "[12.682596210489125, 15.93901395879758]"
$ws.Range("W2").Value = 2.693893893893915
$ws.Range("X2").Value = 0.4644644644644664
$ws.Range("Y2").Value = 4.923323323323363

# --- Row 3 updates ---
$ws.Range("E3").Value = 23.2900000000002
$ws.Range("G3").Value = 0.0004875996742267352
$ws.Range("H3").Value = 0.003059608378842459
$ws.Range("I3").Value = ""
$ws.Range("K3").Value = 4.17044212836067
$ws.Range("L3").Value = "[1.6853379721923272, 6.655546284529013]"
$ws.Range("M3").Value = 0.001048225912857692
$ws.Range("N3").Value = 0.002096451825715384
$ws.Range("O3").Value = -0.1509473947670772
$ws.Range("P3").Value = "[-0.8679475199106932, 0.5660527303765388]"
$ws.Range("Q3").Value = 0.6792990666695284
$ws.Range("R3").Value = 0.6792990666695284
$ws.Range("S3").Value = 12.88830651642785
$ws.Range("T3").Value = "[11.414893121910922, 14.361719910944785]"
$ws.Range("W3").Value = 0.5595195195195259
$ws.Range("X3").Value = -2.098198198198213
$ws.Range("Y3").Value = 3.217237237237265
